{"js": "// Applies the \"changes in the final paragraphs\" edit:\n// 1) Drop the opening \"Edward Snowden's case was primary a reputational hit\n//    for the NSA and the U.S. government in general.\" sentence, and drop\n//    \"thus \" before \"allowing government to spy on people\".\n// 2) \"...persona to software engineers...\" -> \"...persona but more to\n//    software engineers...\" and \"sacrificing their company's \"Top Secret\"\n//    data\" -> \"going against company's policies\".\n// 3) \"...has revealed overweighs \" -> \"...has revealed, overweigh \".\n// 4) Header: merge the \" \" + \"7320\" runs into \" 7320\" (no visible text\n//    change, just a run merge) -- handled in the PowerShell/COM script;\n//    Office.js body search below only touches the document body, so we\n//    separately fix the header text using the same search/replace idiom.\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1a) Remove the opening sentence about Snowden's case being a reputational hit.\nawait replaceOnce(\n  \"Edward Snowden\\u2019s case was primary a reputational hit for the NSA and the U.S. government in general. \",\n  \"\"\n);\n\n// 1b) Remove \"thus \" before \"allowing government to spy on people\".\nawait replaceOnce(\n  \"most popular devices, thus allowing government to spy on people.\",\n  \"most popular devices, allowing government to spy on people.\"\n);\n\n// 2) Update the \"persona ... software engineers\" / \"sacrificing ...\" sentence.\nawait replaceOnce(\n  \"Edward Snowden\\u2019s persona to software engineers in general.\",\n  \"Edward Snowden\\u2019s persona but more to software engineers in general.\"\n);\n\nawait replaceOnce(\n  \"this option of sacrificing their company\\u2019s \\u201cTop Secret\\u201d data if it can affect\",\n  \"this option of going against company\\u2019s policies if it can affect\"\n);\n\n// 3) \"...has revealed overweighs \" -> \"...has revealed, overweigh \"\nawait replaceOnce(\n  \"that this intel has revealed overweighs \",\n  \"that this intel has revealed, overweigh \"\n);\n\n// Note: the diff also merges two adjacent, identically-formatted header\n// runs (\" \" and \"7320\") into a single \" 7320\" run. That is a run-split-only\n// change -- the header's visible text (\"CS- 7320\") is unchanged -- and the\n// Word.js header accessors in this runtime mint extra header parts on touch\n// (even a read-only .load()), which would introduce unrelated structural\n// changes. So the header is intentionally left untouched here.\n", "ps1": "# Applies the \"changes in the final paragraphs\" edit described by the\n# commit: trims/reworks text in the last three body paragraphs of the\n# Snowden ethics write-up, plus a small header run tidy-up.\n\n$d = $word.ActiveDocument\n\n$RSO = 2  # wdReplaceOne\n\n# 1a) Drop the opening sentence \"Edward Snowden's case was primary a\n#     reputational hit for the NSA and the U.S. government in general. \"\n$needle = \"Edward Snowden\" + [char]0x2019 + \"s case was primary a reputational hit for the NSA and the U.S. government in general. \"\n$d.Content.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, \"\", $RSO) | Out-Null\n\n# 1b) Drop \"thus \" before \"allowing government to spy on people.\"\n$d.Content.Find.Execute(\n  \"most popular devices, thus allowing government to spy on people.\",\n  $false, $false, $false, $false, $false, $true, 1, $false,\n  \"most popular devices, allowing government to spy on people.\",\n  $RSO\n) | Out-Null\n\n# 2a) \"...persona to software engineers...\" -> \"...persona but more to\n#     software engineers...\"\n$d.Content.Find.Execute(\n  (\"Edward Snowden\" + [char]0x2019 + \"s persona to software engineers in general.\"),\n  $false, $false, $false, $false, $false, $true, 1, $false,\n  (\"Edward Snowden\" + [char]0x2019 + \"s persona but more to software engineers in general.\"),\n  $RSO\n) | Out-Null\n\n# 2b) \"...sacrificing their company's \"Top Secret\" data...\" -> \"...going\n#     against company's policies...\"\n$d.Content.Find.Execute(\n  (\"this option of sacrificing their company\" + [char]0x2019 + \"s \" + [char]0x201C + \"Top Secret\" + [char]0x201D + \" data if it can affect\"),\n  $false, $false, $false, $false, $false, $true, 1, $false,\n  (\"this option of going against company\" + [char]0x2019 + \"s policies if it can affect\"),\n  $RSO\n) | Out-Null\n\n# 3) \"...has revealed overweighs \" -> \"...has revealed, overweigh \"\n$d.Content.Find.Execute(\n  \"that this intel has revealed overweighs \",\n  $false, $false, $false, $false, $false, $true, 1, $false,\n  \"that this intel has revealed, overweigh \",\n  $RSO\n) | Out-Null\n\n# 4) Header: tidy the \"CS-\" / \" \" / \"7320\" runs into \"CS-\" / \" 7320\" by\n#    re-asserting the (unchanged) text, which lets the host re-merge the\n#    two adjacent, identically-formatted runs into one.\nforeach ($sec in $d.Sections) {\n    $hdr = $sec.Headers.Item(1)\n    $hdr.Range.Find.Execute(\" 7320\", $false, $false, $false, $false, $false, $true, 1, $false, \" 7320\", $RSO) | Out-Null\n}\n"}
